# Daily cryptos data refresh (GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '40.090.30'
$ws.Range('E2').Value = '  -3.05%  '
$ws.Range('D3').Value = '2.339.18'
$ws.Range('E3').Value = '  -4.14%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.81'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '85.49'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -4.42%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.530'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -2.18%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.486'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -2.36%  '
$ws.Range('E10').Value = '  -1.92%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '30.18'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -6.12%  '
$ws.Range('E12').Value = '  +1.18%  '
$ws.Range('D13').Value = '2.698.08'
$ws.Range('E13').Value = '  -4.12%  '
$ws.Range('E14').Value = '  -4.05%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.81'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -4.48%  '
$ws.Range('D16').Value = '2.366.18'
$ws.Range('E16').Value = '  -2.80%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.760'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.80%  '
$ws.Range('D18').Value = '40.065.74'
$ws.Range('E18').Value = '  -2.91%  '
$ws.Range('D19').Value = '0.0₃0904'
$ws.Range('E19').Value = '  -1.95%  '
$ws.Range('E20').Value = '  -1.89%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '68.06'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -5.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.69'
$ws.Range('D22').ClearFormats()
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.46'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('E24').Value = '  -5.01%  '
$ws.Range('E25').Value = '  +0.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.84'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -2.45%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.46'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -2.47%  '
$ws.Range('E28').Value = '  -4.18%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.31'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -2.56%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.96'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.76%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '153.32'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -2.72%  '
$ws.Range('E32').Value = '  -0.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.13'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.79%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.43'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -3.71%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0720'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -3.39%  '
$ws.Range('E36').Value = '  -0.54%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.80'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -3.86%  '
$ws.Range('B38').Value = 'Celestia'
$ws.Range('C38').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '15.82'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -4.61%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0992'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.44%  '
$ws.Range('E40').Value = '  -2.53%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.89'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D42').Value = '1.954.11'
$ws.Range('E42').Value = '  -1.61%  '
$ws.Range('E43').Value = '  -4.62%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0264'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -3.96%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '17.61'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -3.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.50'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.20%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.72'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -5.47%  '
$ws.Range('D48').Value = '2.557.82'
$ws.Range('E48').Value = '  -4.25%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '92.92'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.38%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '70.89'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -3.13%  '
$ws.Range('B51').Value = 'ordi'
$ws.Range('C51').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '63.80'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.17%  '
